$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert a new row above row 3 (shifts old row 3 "stand" down to row 4, etc.)
$ws.Rows.Item(3).Insert()

# New row 3: TreeID field (type=string, name=TreeID, label=TreeID)
$ws.Range("C3").Value = "string"
$ws.Range("E3").Value = "TreeID"
$ws.Range("F3").Value = "TreeID"

# Row 4 (previously row 3, "stand") changes type from integer to string
$ws.Range("C4").Value = "string"
